$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 140 -- shifts existing rows 140:237 down to 141:238
# and grows the sheet dimension from A1:T237 to A1:T238, matching the diff.
$ws.Rows("140:140").Insert()

# The inserted row is blank except for the date-format style that Excel
# carries down from the row below onto column D. Populate the rest of the
# "constant" columns for this product block by copying them from the row
# directly below (the data that used to be row 140 before the insert), then
# overwrite the columns whose values actually changed with the new figures.
$copyCols = @("A","B","C","E","F","G","H","I","J","K","L","Q","R")
foreach ($col in $copyCols) {
    $ws.Range("$col`140").Value2 = $ws.Range("$col`141").Value2
}

$ws.Range("D140").Value2 = 44719
$ws.Range("M140").Value2 = 125
$ws.Range("N140").Value2 = 12000
$ws.Range("O140").Value2 = 12000
$ws.Range("P140").Value2 = 12000
$ws.Range("S140").Value2 = 800
$ws.Range("T140").Value2 = 15
